$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: ASPM, grade B, 2019 - 2023
$ws.Range("A5").Value = "Retaruke at Whanganui Confluence"
$ws.Range("B5").Value = "ASPM"
$ws.Range("C5").Value = "B"
$ws.Range("D5").Value = "2019 - 2023"
$ws.Range("E5").Value = "RepSite"
$ws.Range("F5").Value = 0.41
$ws.Range("G5").Value = 0.4304
$ws.Range("H5").Value = 0.532
$ws.Range("I5").Value = 0.532
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = 0.4015
$ws.Range("M5").Value = 0.51905
$ws.Range("N5").Value = 0.532
$ws.Range("O5").Value = 1779627.51
$ws.Range("P5").Value = 5669030.88
$ws.Range("Q5").Value = "Ruapehu District"
$ws.Range("R5").Value = "Whanganui"
$ws.Range("S5").Value = "Middle Whanganui"
$ws.Range("T5").Value = "Whai_4d"
$ws.Range("U5").Value = ""

# Row 6: MCI, grade C, 2019 - 2023
$ws.Range("A6").Value = "Retaruke at Whanganui Confluence"
$ws.Range("B6").Value = "MCI"
$ws.Range("C6").Value = "C"
$ws.Range("D6").Value = "2019 - 2023"
$ws.Range("E6").Value = "RepSite"
$ws.Range("F6").Value = 111.2
$ws.Range("G6").Value = 108.582
$ws.Range("H6").Value = 113
$ws.Range("I6").Value = 113
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = 107.265
$ws.Range("M6").Value = 112.783
$ws.Range("N6").Value = 113
$ws.Range("O6").Value = 1779627.51
$ws.Range("P6").Value = 5669030.88
$ws.Range("Q6").Value = "Ruapehu District"
$ws.Range("R6").Value = "Whanganui"
$ws.Range("S6").Value = "Middle Whanganui"
$ws.Range("T6").Value = "Whai_4d"
$ws.Range("U6").Value = ""

# Row 7: QMCI, grade C, 2019 - 2023
$ws.Range("A7").Value = "Retaruke at Whanganui Confluence"
$ws.Range("B7").Value = "QMCI"
$ws.Range("C7").Value = "C"
$ws.Range("D7").Value = "2019 - 2023"
$ws.Range("E7").Value = "RepSite"
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 5.165
$ws.Range("H7").Value = 6.418
$ws.Range("I7").Value = 6.418
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = 4.99
$ws.Range("M7").Value = 6.13065
$ws.Range("N7").Value = 6.418
$ws.Range("O7").Value = 1779627.51
$ws.Range("P7").Value = 5669030.88
$ws.Range("Q7").Value = "Ruapehu District"
$ws.Range("R7").Value = "Whanganui"
$ws.Range("S7").Value = "Middle Whanganui"
$ws.Range("T7").Value = "Whai_4d"
$ws.Range("U7").Value = ""
